# Actualización automática de tasas-transfi.xlsx
$wb = $excel.ActiveWorkbook

# --- Update "Hoja1" conversion summary text (cell A1) ---
$wsHoja1 = $wb.Worksheets.Item("Hoja1")
$oldText = $wsHoja1.Range("A1").Text
$newText = $oldText.Replace("1000 Bs = 7.54 = 30708.11 pesos", "1000 Bs = 7.69 = 31201.54 pesos")
$newText = $newText.Replace("30708.11 pesos = 7.54 = 960.38 Bs", "31201.54 pesos = 7.63 = 936.53 Bs")
$wsHoja1.Range("A1").Value = $newText

# --- Update "tasas" sheet rate figures ---
$wsTasas = $wb.Worksheets.Item("tasas")
$wsTasas.Range("N10").Value = 130
$wsTasas.Range("O10").Value = 4056.2
$wsTasas.Range("N12").Value = 4087.54
$wsTasas.Range("O12").Value = 122.69
